# "run all inputs except financial flows"
#
# Adds two new pattern rows to the "Input" sheet's "financial flows" block:
#   - "repay"  (LEMMA pattern), inserted right after the existing "pay" row
#   - "budget" (LEMMA pattern), appended as the new last row of the block
#
# Inserting "repay" pushes the existing "credit"/"debt"/"expenditure"/"income"
# rows down by one; "budget" is then added as a brand new row after "income".
# The single-cell-range list validation on the pattern columns is extended to
# cover the two extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# --- insert "repay" right after "pay" (row 17 -> new row 18) ---------------
$ws.Rows.Item(18).Insert()
$ws.Range("B18").Value = "FINANCIAL_FLOW"
$ws.Range("D18").Value = "repay"
$ws.Range("E18").Value = "LEMMA"
$ws.Range("F18").Value = "repay"

# --- insert a fresh row at 22 (inside the still-valid 3:22 validation area) -
# This pushes the "income" row (now at 22) down to 23, extends the single
# list-validation range from row 22 to row 23, and leaves row 22 blank for
# us to fill back in.
$ws.Rows.Item(22).Insert()

$ws.Range("B22").Value = "FINANCIAL_FLOW"
$ws.Range("D22").Value = "income"
$ws.Range("E22").Value = "LEMMA"
$ws.Range("F22").Value = "income"

$ws.Range("B23").Value = "FINANCIAL_FLOW"
$ws.Range("D23").Value = "budget"
$ws.Range("E23").Value = "LEMMA"
$ws.Range("F23").Value = "budget"
